$d = $word.ActiveDocument

# 1. Remove the stale "_GoBack" bookmark (Word drops/relocates this automatically once the
#    document is edited somewhere else; the old marker at the "V1.972" line is gone afterwards).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Proofread fix: "By default a BASTION zone" is missing the comma after the introductory
#    adverb -> "By default, a BASTION zone".
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("By default a BASTION zone", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "By default, a BASTION zone", 2)
